# edit.ps1 -- applies the "TimesNewToman" Word-doc diff via Word COM interop.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Unveiling the Enigmatic Cosmos: A Journey of Exploration",
    $true, $false, $false, $false, $false, $true, 1, $false,
    'Tracing the Footprints of Chemistry: Unveiling the World of Matter and Transformation', 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Author name -> "Dr. Isabel Goddard"
# ---------------------------------------------------------------------------
$authorPara = $d.Paragraphs.Item(2).Range
$authorPara.Text = "Dr. Isabel Goddard"

# ---------------------------------------------------------------------------
# 3) Email -> "goddard_isa@edu.net"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "isabelle-dubois@astronomy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "goddard_isa@edu", 2) | Out-Null
$d.Content.Find.Execute(
    "org",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "net", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Main body paragraph (paragraph 5) - full rewrite
# ---------------------------------------------------------------------------
$para5Text = 'In the realm of science, chemistry stands out as a captivating discipline that delves into the intricate world of matter and its transformations. Chemistry is the study of substances, their properties, and the changes they undergo. It seeks to unravel the mysteries of how and why matter behaves the way it does, revealing the fundamental principles that govern the world around us.' + [char]11 + '' + [char]11 + 'From the smallest atoms to the vast universe, chemistry plays a pivotal role in shaping our understanding of the natural world. It underpins our comprehension of materials, their interactions, and their applications in diverse fields ranging from medicine and engineering to agriculture and environmental science. Chemistry is a field where human ingenuity and scientific inquiry converge to unlock the secrets of matter, propelling technological advancements and solving global challenges.' + [char]11 + '' + [char]11 + 'Whether it''s the intricate dance of electrons within molecules or the mesmerizing transformations that occur during chemical reactions, chemistry offers a window into the hidden realm of matter. It challenges our curiosity, igniting a desire to explore the unseen and unravel the mysteries that lie beneath the surface.' + [char]11 + '' + [char]11 + 'Introduction Continued:' + [char]11 + '' + [char]11 + 'Chemistry has been an integral part of human civilization since ancient times, with alchemists seeking to transmute elements and create the elusive philosopher''s stone. Over the centuries, the understanding of chemistry has evolved dramatically, thanks to the contributions of brilliant minds like Antoine Lavoisier, Marie Curie, and Linus Pauling. These pioneers laid the groundwork for modern chemistry, establishing fundamental principles and developing analytical techniques that have revolutionized our knowledge of matter.' + [char]11 + '' + [char]11 + 'Today, chemistry continues to be a dynamic and evolving field, with new discoveriesBu Duan  being made. From the development of life-saving drugs to the creation of sustainable materials, chemistry is at the forefront of scientific innovation. As we delve deeper into the intricacies of matter and its interactions, we unlock new possibilities for addressing pressing global issues such as climate change, energy security, and food production.' + [char]11 + '' + [char]11 + 'Introduction Continued:' + [char]11 + '' + [char]11 + 'The study of chemistry is not merely about memorizing formulas and equations; it''s about fostering a mindset of critical thinking, problem-solving, and creative exploration. Chemistry encourages students to ask questions, investigate phenomena, and develop logical reasoning skills. It cultivates an appreciation for the elegance and beauty of the natural world, while also highlighting the profound impact that chemistry has on our lives and the environment.'
$body5 = $d.Paragraphs.Item(5).Range
$body5.Text = $para5Text

# ---------------------------------------------------------------------------
# 5) Summary body (paragraph 7) - full rewrite
# ---------------------------------------------------------------------------
$summaryPara = $d.Paragraphs.Item(7).Range
$summaryPara.Text = 'In this essay, we explored the captivating world of chemistry, unveiling its fundamental principles and highlighting its significance in shaping our understanding of matter and its transformations. Chemistry is a field that challenges our curiosity, ignites our imagination, and empowers us to address global challenges. By delving into the intricacies of matter, we unlock new possibilities for innovation, sustainability, and human progress.'

# ---------------------------------------------------------------------------
# 6) Add a trailing empty paragraph at the end of the document
# ---------------------------------------------------------------------------
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 7) Fix the font across the whole document: TimesNewToman -> Times New Roman
# ---------------------------------------------------------------------------
$full = $d.Range(0, $d.Content.End)
$full.Font.Name = "Times New Roman"
